# Generate Report for Handoff
# Refresh the handoff UUID (585d9800-... -> bef8f109-...) and the
# associated handoff timestamps across the Overview / zh-cn / de-de sheets,
# keeping each hyperlink's underlying target URL (and therefore its
# relationship id) unchanged -- only the visible text / display changes.

$wb = $excel.ActiveWorkbook

$oldGuid = "585d9800-508d-44bb-81ca-85145cb74480"
$newGuid = "bef8f109-06de-4819-9d2b-9e3f49d29f6c"

$oldMd  = "$oldGuid.md"
$newMd  = "$newGuid.md"

$oldZh  = "$oldGuid.fe1d474cfb6247d8f3ef93a3c7321f20bc2f3c6d.zh-cn.xlf"
$newZh  = "$newGuid.fa4a2a1344bd61e1b34cb91c9d4bcf5cd9de2b7f.zh-cn.xlf"

$oldDe  = "$oldGuid.fe1d474cfb6247d8f3ef93a3c7321f20bc2f3c6d.de-de.xlf"
$newDe  = "$newGuid.fa4a2a1344bd61e1b34cb91c9d4bcf5cd9de2b7f.de-de.xlf"

$newHandoffDate = "2016-03-25 01:21:02"
$newZhDate      = "2016-03-25 01:20:57"

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/80e36df775d6e846e0e1cca3d695f8614f9f5137/e2e/$oldMd"

$wsOverview.Range("A2").Value2 = $newMd
$wsOverview.Range("D2").Value2 = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMd) | Out-Null

# ---- zh-cn sheet -------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhMdAddress  = "https://github.com/OpenLocalizationTest/oltest/blob/80e36df775d6e846e0e1cca3d695f8614f9f5137/e2e/$oldMd"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b9fc44427d434901618968d99c66621d6238f6d1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZh"

$wsZh.Range("A2").Value2 = $newMd
$wsZh.Range("D2").Value2 = $newZh
$wsZh.Range("E2").Value2 = $newZhDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhMdAddress, "", "", $newMd) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfAddress, "", "", $newZh) | Out-Null

# ---- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deMdAddress  = "https://github.com/OpenLocalizationTest/oltest/blob/80e36df775d6e846e0e1cca3d695f8614f9f5137/e2e/$oldMd"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/de81313fa482d981e94551accab49e9d7fdc98c2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDe"

$wsDe.Range("A2").Value2 = $newMd
$wsDe.Range("D2").Value2 = $newDe
$wsDe.Range("E2").Value2 = $newHandoffDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deMdAddress, "", "", $newMd) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfAddress, "", "", $newDe) | Out-Null

Write-Host "Report generated for handoff."
